# Insert a new daily data point (2026/01/26, 14:00, rank 14/156) into the
# Sei1 ranking sheet. The new row is inserted at row 727 (between the
# existing 2026/01/26 20:00 entry at row 726 and the 2026/12/29 entry that
# used to be row 727), pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 727, shifting 727:768 down
# to 728:769 (and updating the sheet dimension from D768 to D769).
$ws.Range("A727").EntireRow.Insert()

# Populate the newly inserted row. Force column A to be stored as plain
# text (matching the rest of the date column) instead of being
# auto-converted to a date serial value.
$ws.Range("A727").NumberFormat = "@"
$ws.Range("A727").Value = "2026/01/26"
$ws.Range("B727").Value = "月"
$ws.Range("C727").Value = 14
$ws.Range("D727").Value = 156
